# Update "Response sizes" Sheet1 statistics and selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New Time [ms] (column B) and Size [MB] (column C) values for rows 2-21.
$values = @(
    @(609, 36),
    @(1013, 73),
    @(1510, 109),
    @(2022, 146),
    @(2514, 182),
    @(3057, 219),
    @(3507, 255),
    @(4013, 292),
    @(4591, 328),
    @(5038, 364),
    @(5740, 401),
    @(6073, 438),
    @(6680, 474),
    @(7361, 511),
    @(7727, 547),
    @(8203, 583),
    @(8986, 620),
    @(9519, 657),
    @(9685, 693),
    @(10230, 729)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
}

# Update the active cell / selection shown in the saved sheet view.
$ws.Activate()
$ws.Range("H18").Select()
